$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")
$ws.Activate()

# Root-cause change: the upfront payment in Q7 (year 2027) was removed.
# All downstream formulas (Q12, Q16:Y20, Q22:Y22, AB26:AB28, Q29:Y29, etc.)
# are formula-driven and recalculate automatically once this input changes.
$ws.Range("Q7").Value = 0

# Update the on-screen selection to match the saved view state
# (bottom-right frozen pane: U21:AA22, with U21 as the active cell).
$ws.Range("U21:AA22").Select()
